$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10482567.6771324
$ws.Range("C2").Value = 63838446353.3055
$ws.Range("D2").Value = 8394316.75268444

$ws.Range("B3").Value = 1619327.16855146
$ws.Range("C3").Value = 1667833.87656769
$ws.Range("D3").Value = 1622476.73057767

$wb.Save()
